# Time Log.xlsx - add the missing "Testing" time entry for 2014-09-09
# (row 37 of Sheet1, which was previously a blank placeholder row), and
# move the active-cell selection down to C38 to match where the user
# left off after filling the row in.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Date (stored as an Excel serial date number, same as the other rows)
$ws.Range("A37").Value = 41891
# Start Time / Stop Time (fractions of a day)
$ws.Range("B37").Value = 0.90555555555555556
$ws.Range("C37").Value = 0.93819444444444444
# Interruption (minutes)
$ws.Range("D37").Value = 5
# Delta - same shared formula as the rest of column E
$ws.Range("E37").Formula = "=IF(AND(NOT(ISBLANK(B37)),NOT(ISBLANK(C37))), (C37-B37) * 24 - D37/60, """")"
# Activity
$ws.Range("F37").Value = "Testing"

# Leave the selection where the user would be after entering this row
$ws.Range("C38").Select()
